$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Location'
$ws.Range("B1").Value = 'Site'
$ws.Range("C1").Value = 'Exposure period'
$ws.Range("D1").Value = 'Notes'

$ws.Range("A2").Value = 'Bentleigh'
$ws.Range("B2").Value = 'Il Centro Deli  5/284/292 Centre Rd, Bentleigh VIC 3204'
$ws.Range("C2").Value = '22/12/20 12:00pm-12:30pm'
$ws.Range("D2").Value = 'Case shopped in store'

$ws.Range("A3").Value = 'Black Rock'
$ws.Range("B3").Value = 'Woolworths Metro  40 Bluff Road, Black Rock VIC 3193'
$ws.Range("C3").Value = '30/12/20 5:30pm-5:55pm'
$ws.Range("D3").Value = 'Case shopped'

$ws.Range("A4").Value = 'Box Hill South'
$ws.Range("B4").Value = 'Bunnings  259 Middleborough Road, Box Hill South VIC 3128'
$ws.Range("C4").Value = '30/12/20 12:00pm-12:40pm'
$ws.Range("D4").Value = 'Case shopped'

$ws.Range("A5").Value = 'Brighton'
$ws.Range("B5").Value = 'Brighton Beach  Brighton, VIC 3186'
$ws.Range("C5").Value = '26/12/20 12:00pm-1:00pm'
$ws.Range("D5").Value = 'Case attended beach'

$ws.Range("A6").Value = 'Brighton'
$ws.Range("B6").Value = 'Brighton Beach  Brighton, VIC 3186'
$ws.Range("C6").Value = '29/12/20 12:00pm-3:00pm'
$ws.Range("D6").Value = 'Case attended beach'

$ws.Range("A7").Value = 'Burwood East'
$ws.Range("B7").Value = 'Coles Burwood, Burwood Highway & Blackburn Road'
$ws.Range("C7").Value = '28/12/20 6.30pm - 7pm'
$ws.Range("D7").Value = 'Case shopped in store'

$ws.Range("A8").Value = 'Burwood East'
$ws.Range("B8").Value = 'Kmart Burwood, 172 Burwood Highway'
$ws.Range("C8").Value = '28/12/20 6.15pm - 6.30pm'
$ws.Range("D8").Value = 'Case shopped in store'

$ws.Range("A9").Value = 'Camberwell'
$ws.Range("B9").Value = 'Coles Middle Camberwell, 751 Riversdale Road'
$ws.Range("C9").Value = '28/12/20 12pm - 12.30pm'
$ws.Range("D9").Value = 'Case shopped in store'

$ws.Range("A10").Value = 'Camberwell'
$ws.Range("B10").Value = 'Fu Lin Asian Grocery Supermarket  1397 Toorak Road, Camberwell VIC 3124'
$ws.Range("C10").Value = '30/12/20 2:30pm-2:45pm'
$ws.Range("D10").Value = 'Case shopped'

$ws.Range("A11").Value = 'Cape Schank'
$ws.Range("B11").Value = 'National Golf Club  The Cups Drive, Cape Schanck VIC 3939'
$ws.Range("C11").Value = '30/12/20 11.40am-1.40pm'
$ws.Range("D11").Value = 'Case attended course'

$ws.Range("A12").Value = 'Cheltenham'
$ws.Range("B12").Value = 'Aldi Cheltenham  280/282 Bay Road, Cheltenham VIC 3192'
$ws.Range("C12").Value = '29/12/20 2:00pm-2:30pm'
$ws.Range("D12").Value = 'Case shopped in store'

$ws.Range("A13").Value = 'Cheltenham'
$ws.Range("B13").Value = 'Angus and Cootes Jeweller  Southland Shopping Centre, Shop 2096/1239, Nepean Hwy, Cheltenham VIC 3192'
$ws.Range("C13").Value = '28/12/2020 2:30pm-2:50pm'
$ws.Range("D13").Value = 'Case shopped in store'

$ws.Range("A14").Value = 'Cheltenham'
$ws.Range("B14").Value = 'Bodero Southland Shopping Centre, 1239 Nepean Hwy'
$ws.Range("C14").Value = '22/12/20 6.45pm - 7pm'
$ws.Range("D14").Value = 'Case shopped in store'

$ws.Range("A15").Value = 'Cheltenham'
$ws.Range("B15").Value = 'Chemist Warehouse Cheltenham, 326/330 Charman Rd'
$ws.Range("C15").Value = '03/01/21, 3.30pm - 3.45pm'
$ws.Range("D15").Value = 'Case shopped in store'

$ws.Range("A16").Value = 'Cheltenham'
$ws.Range("B16").Value = 'Coles, Westfield Southland'
$ws.Range("C16").Value = '22/12/20 11:50am-12:10pm'
$ws.Range("D16").Value = 'Case shopped in store'

$ws.Range("A17").Value = 'Cheltenham'
$ws.Range("B17").Value = 'Cotton On, Southland Shopping Centre 1239 Nepean Hwy'
$ws.Range("C17").Value = '22/12/20 12.15pm - 12.45pm'
$ws.Range("D17").Value = 'Case visited venue'

$ws.Range("A18").Value = 'Cheltenham'
$ws.Range("B18").Value = 'Honey Birdette Southland  Shop 2209/1239, Southland Shopping Centre, Cheltenham VIC 3192'
$ws.Range("C18").Value = '22/12/2020 3:50pm-4:05pm'
$ws.Range("D18").Value = 'Case shopped in store'

$ws.Range("A19").Value = 'Cheltenham'
$ws.Range("B19").Value = 'Kmart Southland Shopping Centre, 1239 Nepean Highway'
$ws.Range("C19").Value = '22/12/20 6.30pm - 6.45pm'
$ws.Range("D19").Value = 'Case shopped in store'

$ws.Range("A20").Value = 'Cheltenham'
$ws.Range("B20").Value = 'Kmart Southland Shopping Centre, 1239 Nepean Highway'
$ws.Range("C20").Value = '28/12/20 2.30pm-3pm'
$ws.Range("D20").Value = 'Case shopped in store'

$ws.Range("A21").Value = 'Cheltenham'
$ws.Range("B21").Value = 'Mecca Southland  Shop 2011/2013, Southland Shopping Centre, Cheltenham VIC 3192'
$ws.Range("C21").Value = '22/12/2020 3:30pm-3:50pm'
$ws.Range("D21").Value = 'Case shopped in store'

$ws.Range("A22").Value = 'Cheltenham'
$ws.Range("B22").Value = 'Myer, Southland Shopping Centre 1239 Nepean Hwy'
$ws.Range("C22").Value = '22/12/20 10.30am - 11am'
$ws.Range("D22").Value = 'Case visited venue'

$ws.Range("A23").Value = 'Cheltenham'
$ws.Range("B23").Value = 'Specsavers, 1004-1005 Westfield Southland'
$ws.Range("C23").Value = '22/12/20 11:00am-1145am'
$ws.Range("D23").Value = 'Case shopped in store'

$ws.Range("A24").Value = 'Clayton'
$ws.Range("B24").Value = 'Kmart - M-City, 2107 Dandenong Road, Clayton'
$ws.Range("C24").Value = '30/12/20 7pm - 7.30pm'
$ws.Range("D24").Value = 'Case shopped at store'

$ws.Range("A25").Value = 'Clayton'
$ws.Range("B25").Value = 'Woolworths - M-City, 2107 Dandenong Road, Clayton'
$ws.Range("C25").Value = '30/12/20 7.30pm - 745pm'
$ws.Range("D25").Value = 'Case shopped at store'

$ws.Range("A26").Value = 'Forest Hill'
$ws.Range("B26").Value = 'Forest Hill Chase Shopping Centre 270 Canterbury Rd, Forest Hill VIC 3131'
$ws.Range("C26").Value = '28/12/20 12:00pm-2:00pm'
$ws.Range("D26").Value = '1210hrs Food court 30min; 1250hrs TKMaxx 15min; 1310hrs Target 20min; 1340hrs Woolworths 15min'

$ws.Range("A27").Value = 'Fountain Gate Shopping Centre'
$ws.Range("B27").Value = 'Kmart, Big W, Target, Millers, King of Gifts, Lo Costa  25-55 Overland Drive, Narre Warren VIC 3805'
$ws.Range("C27").Value = '26/12/20 9:00am-11:00am'
$ws.Range("D27").Value = $null

$ws.Range("A28").Value = 'Frankston'
$ws.Range("B28").Value = 'TK Maxx Frankston, 10 Shannon Street, Bayside Shopping Centre'
$ws.Range("C28").Value = '31/12/20 2pm - 3pm'
$ws.Range("D28").Value = 'Case shopped at store'

$ws.Range("A29").Value = 'Glen Waverley'
$ws.Range("B29").Value = 'Mocha Jos  87 Kingsway, Glen Waverley VIC 3150'
$ws.Range("C29").Value = '28/12/20 1:30pm-1:45pm'
$ws.Range("D29").Value = $null

$ws.Range("A30").Value = 'Hallam'
$ws.Range("B30").Value = 'Coles Hallam  2 Princes Domain Drive, Hallam VIC 3803'
$ws.Range("C30").Value = '30/12/20 6:15am-6:30am'
$ws.Range("D30").Value = 'Case shopped in store'

$ws.Range("A31").Value = 'Lakes Entrance'
$ws.Range("B31").Value = 'Blue Riviera Hire Boats  Marine Parade, Lakes Entrance VIC 3909'
$ws.Range("C31").Value = '29/12/20 11:15am-12:15pm'
$ws.Range("D31").Value = 'Case hired a boat'

$ws.Range("A32").Value = 'Lakes Entrance'
$ws.Range("B32").Value = 'Central Hotel Lakes Entrance  321 Esplanade, Lakes Entrance VIC 3909'
$ws.Range("C32").Value = '30/12/20 5:00pm-6:30pm'
$ws.Range("D32").Value = 'Case attended outside premises'

$ws.Range("A33").Value = 'Lakes Entrance'
$ws.Range("B33").Value = 'Darcey Annas Beach Cafe Kiosk Gift Shop Gallery  426 Main Beach Walk Surf Life Saving, Lakes Entrance VIC 3909'
$ws.Range("C33").Value = '30/12/20 11:15am-11:20am'
$ws.Range("D33").Value = 'Case picked up takeaway'

$ws.Range("A34").Value = 'Lakes Entrance'
$ws.Range("B34").Value = 'Woolworths Lakes Entrance 371 Esplanade, Lakes Entrance VIC 3909'
$ws.Range("C34").Value = '30/12/20 6:00pm-6:15pm'
$ws.Range("D34").Value = 'Case shopped in store'

$ws.Range("A35").Value = 'Melbourne'
$ws.Range("B35").Value = 'Federation Square  Swanston & Flinders Streets, Melbourne VIC 3000'
$ws.Range("C35").Value = '23/12/2- 11:00pm-11:30pm'
$ws.Range("D35").Value = 'Case attended Federation Square'

$ws.Range("A36").Value = 'Mentone'
$ws.Range("B36").Value = 'Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194'
$ws.Range("C36").Value = '23/12/20 11:00am-11:30am'
$ws.Range("D36").Value = 'Case shopped in store'

$ws.Range("A37").Value = 'Mentone'
$ws.Range("B37").Value = 'Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194'
$ws.Range("C37").Value = '29/12/20 07:30am-08:00am'
$ws.Range("D37").Value = 'Case shopped in store'

$ws.Range("A38").Value = 'Mentone'
$ws.Range("B38").Value = 'Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194'
$ws.Range("C38").Value = '31/12/20 08:00am-08:30am'
$ws.Range("D38").Value = 'Case shopped in store'

$ws.Range("A39").Value = 'Mentone'
$ws.Range("B39").Value = 'Mentone/Parkdale Beach'
$ws.Range("C39").Value = '27/12/20 10:00am-4:30pm'
$ws.Range("D39").Value = $null

$ws.Range("A40").Value = 'Mentone'
$ws.Range("B40").Value = 'Woolworths Mentone  105-111 Balcombe Road, Mentone VIC 3194'
$ws.Range("C40").Value = '23/12/20 2:45pm-3:05pm'
$ws.Range("D40").Value = 'Case shopped in store'

$ws.Range("A41").Value = 'Moorabbin'
$ws.Range("B41").Value = 'COSTCO Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194'
$ws.Range("C41").Value = '30/12/20 10:45am-12:15pm'
$ws.Range("D41").Value = 'Case shopped in store'

$ws.Range("A42").Value = 'Moorabbin'
$ws.Range("B42").Value = 'COSTCO Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194'
$ws.Range("C42").Value = '30/12/20 4:00m- 5:50pm'
$ws.Range("D42").Value = 'Case shopped in store'

$ws.Range("A43").Value = 'Mordialloc'
$ws.Range("B43").Value = 'Woodlands Golf Club  109 White Street Mordialloc VIC 3195'
$ws.Range("C43").Value = '23/12/20 8:00am-2:00pm'
$ws.Range("D43").Value = 'Case attended course'

$ws.Range("A44").Value = 'Mordialloc'
$ws.Range("B44").Value = 'Woodlands Golf Club  109 White Street Mordialloc VIC 3195'
$ws.Range("C44").Value = '28/12/20 12:00pm-6:00pm'
$ws.Range("D44").Value = 'Case attended course'

$ws.Range("A45").Value = 'Mount Martha'
$ws.Range("B45").Value = 'Mount Martha Fine Foods, 34 Lochiel Ave'
$ws.Range("C45").Value = '31/12/20 3pm - 3.15pm'
$ws.Range("D45").Value = 'Takeaway coffee'

$ws.Range("A46").Value = 'Mount Waverley'
$ws.Range("B46").Value = 'Ritchies IGA  283 Stephensons Road, Mount Waverley VIC 3149'
$ws.Range("C46").Value = '30/12/20 2:00pm-2:30pm'
$ws.Range("D46").Value = 'Case shopped for half an hour'

$ws.Range("A47").Value = 'Oakleigh'
$ws.Range("B47").Value = 'Bunnings Oakleigh  1041 Centre Road, Oakleigh South'
$ws.Range("C47").Value = '30/12/20 11:00am-11:30am'
$ws.Range("D47").Value = 'Case shopped for 30 minutes'

$ws.Range("A48").Value = 'Oakleigh'
$ws.Range("B48").Value = 'Katialo restaurant  8 Eaton Mall, Oakleigh VIC 3166'
$ws.Range("C48").Value = '28/12/20 7:00pm-7:10pm'
$ws.Range("D48").Value = $null

$ws.Range("A49").Value = 'Springvale'
$ws.Range("B49").Value = 'IKEA Springvale, 917 Princes Hwy'
$ws.Range("C49").Value = '29/12/20 4pm - 6pm'
$ws.Range("D49").Value = 'Case shopped at store and dined at cafe'

$ws.Range("A50").Value = 'Wonthaggi'
$ws.Range("B50").Value = 'Wonthaggi Plaza Shopping centre  2 Biggs Drive, Wonthaggi VIC 3995'
$ws.Range("C50").Value = '28/12/20 1:30pm-2.30pm'
$ws.Range("D50").Value = 'Kmart- shopped for 15 mins'

